# "new field added in territory import"
#
# A new "CountryName" column is inserted as the first column of the
# territory import template (shifting StateName / DistrictName / CityName /
# IsActive one column to the right, and the IsActive list-validation with
# them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting a whole column before A pushes the existing data, column
# widths, and the data validation range one column to the right for us.
$ws.Range("A1").EntireColumn.Insert()

# Copy the header formatting (fill/style) used by the other header cells
# onto the new A1 header, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "CountryName"

# Match the new column's width from the target layout.
$ws.Range("A1").ColumnWidth = 19.8333333333

# Match the saved selection in the target workbook view.
$ws.Range("D8").Select()
